# Refresh cached Universalis market-price snapshots (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -- columns H:N)
# in the leve-profit table on each job sheet. Pure data refresh from the
# scheduled pricing run -- no formulas or other columns are touched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 7008
$ws.Range("I8").Value = 7008
$ws.Range("K8").Value = 21024
$ws.Range("M8").Value = -20885
# Row 39
$ws.Range("H39").Value = 429.73685
$ws.Range("I39").Value = 88.083336
$ws.Range("J39").Value = 1015.4286
$ws.Range("K39").Value = 264.250008
$ws.Range("L39").Value = 3046.2858
$ws.Range("M39").Value = 31.74999200000002
$ws.Range("N39").Value = -3638.2858
# Row 113
$ws.Range("H113").Value = 5264.4707
$ws.Range("I113").Value = 4355.5557
$ws.Range("J113").Value = 6287
$ws.Range("K113").Value = 4355.5557
$ws.Range("L113").Value = 6287
$ws.Range("M113").Value = -1101.5557
$ws.Range("N113").Value = -12795
# Row 132
$ws.Range("H132").Value = 21370904
$ws.Range("I132").Value = 2416963.5
$ws.Range("J132").Value = 166684430
$ws.Range("K132").Value = 7250890.5
$ws.Range("L132").Value = 500053290
$ws.Range("M132").Value = -7248360.5
$ws.Range("N132").Value = -500058350
# Row 137
$ws.Range("H137").Value = 1381.0625
$ws.Range("I137").Value = 1014.4375
$ws.Range("K137").Value = 3043.3125
$ws.Range("M137").Value = -493.3125

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2548.111
$ws.Range("I2").Value = 2615
$ws.Range("K2").Value = 2615
$ws.Range("M2").Value = -2502
# Row 63
$ws.Range("H63").Value = 3975
$ws.Range("I63").Value = 2950
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2950
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -2264
$ws.Range("N63").Value = -6372
# Row 66
$ws.Range("H66").Value = 3975
$ws.Range("I66").Value = 2950
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 14750
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -11318
$ws.Range("N66").Value = -31864
# Row 74
$ws.Range("H74").Value = 22476.596
$ws.Range("I74").Value = 34376.8
$ws.Range("J74").Value = 1476.2354
$ws.Range("K74").Value = 34376.8
$ws.Range("L74").Value = 1476.2354
$ws.Range("M74").Value = -33502.8
$ws.Range("N74").Value = -3224.2354
# Row 77
$ws.Range("H77").Value = 22476.596
$ws.Range("I77").Value = 34376.8
$ws.Range("J77").Value = 1476.2354
$ws.Range("K77").Value = 171884
$ws.Range("L77").Value = 7381.177
$ws.Range("M77").Value = -167516
$ws.Range("N77").Value = -16117.177
# Row 102
$ws.Range("H102").Value = 1497.56
$ws.Range("I102").Value = 1289.3125
$ws.Range("J102").Value = 1867.7778
$ws.Range("K102").Value = 1289.3125
$ws.Range("L102").Value = 1867.7778
$ws.Range("M102").Value = 332.6875
$ws.Range("N102").Value = -5111.7778
# Row 109
$ws.Range("H109").Value = 29771.428
$ws.Range("J109").Value = 29771.428
$ws.Range("L109").Value = 29771.428
$ws.Range("N109").Value = -32545.428
# Row 116
$ws.Range("H116").Value = 2548.111
$ws.Range("I116").Value = 2615
$ws.Range("K116").Value = 2615
$ws.Range("M116").Value = -321
# Row 132
$ws.Range("H132").Value = 2090.125
$ws.Range("I132").Value = 2098.1853
$ws.Range("J132").Value = 2073.3845
$ws.Range("K132").Value = 6294.5559
$ws.Range("L132").Value = 6220.1535
$ws.Range("M132").Value = -3764.5559
$ws.Range("N132").Value = -11280.1535

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2548.111
$ws.Range("I3").Value = 2615
$ws.Range("K3").Value = 2615
$ws.Range("M3").Value = -2501
# Row 99
$ws.Range("H99").Value = 1779.48
$ws.Range("I99").Value = 1736.875
$ws.Range("J99").Value = 1855.2222
$ws.Range("K99").Value = 1736.875
$ws.Range("L99").Value = 1855.2222
$ws.Range("M99").Value = -238.875
$ws.Range("N99").Value = -4851.2222
# Row 112
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1015
$ws.Range("I16").Value = 830.3333
$ws.Range("J16").Value = 1199.6666
$ws.Range("K16").Value = 830.3333
$ws.Range("L16").Value = 1199.6666
$ws.Range("M16").Value = -543.3333
$ws.Range("N16").Value = -1773.6666
# Row 31
$ws.Range("H31").Value = 16130385
$ws.Range("I31").Value = 21739952
$ws.Range("J31").Value = 2880.625
$ws.Range("K31").Value = 21739952
$ws.Range("L31").Value = 2880.625
$ws.Range("M31").Value = -21739657
$ws.Range("N31").Value = -3470.625
# Row 34
$ws.Range("H34").Value = 16130385
$ws.Range("I34").Value = 21739952
$ws.Range("J34").Value = 2880.625
$ws.Range("K34").Value = 21739952
$ws.Range("L34").Value = 2880.625
$ws.Range("M34").Value = -21739750
$ws.Range("N34").Value = -3284.625
# Row 113
$ws.Range("H113").Value = 1015
$ws.Range("I113").Value = 830.3333
$ws.Range("J113").Value = 1199.6666
$ws.Range("K113").Value = 830.3333
$ws.Range("L113").Value = 1199.6666
$ws.Range("M113").Value = 1339.6667
$ws.Range("N113").Value = -5539.6666
# Row 141
$ws.Range("H141").Value = 66813.5
$ws.Range("J141").Value = 72030.28999999999
$ws.Range("L141").Value = 72030.28999999999
$ws.Range("N141").Value = -82390.28999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 146.6
$ws.Range("I11").Value = 96
$ws.Range("J11").Value = 602
$ws.Range("K11").Value = 288
$ws.Range("L11").Value = 1806
$ws.Range("M11").Value = -148
$ws.Range("N11").Value = -2086
# Row 108
$ws.Range("H108").Value = 2601.889
$ws.Range("I108").Value = 1345.2858
$ws.Range("J108").Value = 7000
$ws.Range("K108").Value = 4035.8574
$ws.Range("L108").Value = 21000
$ws.Range("M108").Value = -1155.8574
$ws.Range("N108").Value = -26760
# Row 109
$ws.Range("H109").Value = 1783.3
$ws.Range("I109").Value = 814.0833
$ws.Range("J109").Value = 3237.125
$ws.Range("K109").Value = 2442.2499
$ws.Range("L109").Value = 9711.375
$ws.Range("M109").Value = -1402.2499
$ws.Range("N109").Value = -11791.375
# Row 131
$ws.Range("H131").Value = 916.66
$ws.Range("I131").Value = 900
$ws.Range("J131").Value = 916.8283
$ws.Range("K131").Value = 2700
$ws.Range("L131").Value = 2750.4849
$ws.Range("M131").Value = 2340
$ws.Range("N131").Value = -12830.4849
# Row 137
$ws.Range("H137").Value = 16511807
$ws.Range("J137").Value = 18296362
$ws.Range("L137").Value = 54889086
$ws.Range("N137").Value = -54899286

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 1669.5
$ws.Range("I9").Value = 1669.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1669.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -1499.5
$ws.Range("N9").ClearContents()
# Row 80
$ws.Range("H80").Value = 2430.5
$ws.Range("I80").Value = 2352.875
$ws.Range("J80").Value = 2585.75
$ws.Range("K80").Value = 2352.875
$ws.Range("L80").Value = 2585.75
$ws.Range("M80").Value = -1354.875
$ws.Range("N80").Value = -4581.75
# Row 83
$ws.Range("H83").Value = 2430.5
$ws.Range("I83").Value = 2352.875
$ws.Range("J83").Value = 2585.75
$ws.Range("K83").Value = 11764.375
$ws.Range("L83").Value = 12928.75
$ws.Range("M83").Value = -6772.375
$ws.Range("N83").Value = -22912.75

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 46249.953
$ws.Range("I22").Value = 200557.6
$ws.Range("J22").Value = 865.35297
$ws.Range("K22").Value = 200557.6
$ws.Range("L22").Value = 865.35297
$ws.Range("M22").Value = -200262.6
$ws.Range("N22").Value = -1455.35297
# Row 27
$ws.Range("H27").Value = 46249.953
$ws.Range("I27").Value = 200557.6
$ws.Range("J27").Value = 865.35297
$ws.Range("K27").Value = 200557.6
$ws.Range("L27").Value = 865.35297
$ws.Range("M27").Value = -200450.6
$ws.Range("N27").Value = -1079.35297
# Row 132
$ws.Range("H132").Value = 2648.9492
$ws.Range("I132").Value = 3116.3157
$ws.Range("J132").Value = 1803.238
$ws.Range("K132").Value = 9348.947100000001
$ws.Range("L132").Value = 5409.714
$ws.Range("M132").Value = -6818.947100000001
$ws.Range("N132").Value = -10469.714

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2888.3572
$ws.Range("I132").Value = 3548.32
$ws.Range("J132").Value = 1917.8235
$ws.Range("K132").Value = 10644.96
$ws.Range("L132").Value = 5753.470499999999
$ws.Range("M132").Value = -8114.960000000001
$ws.Range("N132").Value = -10813.4705
# Row 136
$ws.Range("H136").Value = 878.2125
$ws.Range("I136").Value = 519.46155
$ws.Range("J136").Value = 1544.4642
$ws.Range("K136").Value = 1558.38465
$ws.Range("L136").Value = 4633.392599999999
$ws.Range("M136").Value = 991.61535
$ws.Range("N136").Value = -9733.392599999999
